# Scheduled-runner market data refresh: updates cached market-board figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a handful of leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets, per the latest price pull.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# hunk @ diff line 3834
$ws.Range("H64").Value = 3079.9
$ws.Range("I64").Value = 3079.8
$ws.Range("J64").Value = 3080
$ws.Range("K64").Value = 3079.8
$ws.Range("L64").Value = 3080
$ws.Range("M64").Value = -2831.8
$ws.Range("N64").Value = -3576
# hunk @ diff line 3987
$ws.Range("H67").Value = 3079.9
$ws.Range("I67").Value = 3079.8
$ws.Range("J67").Value = 3080
$ws.Range("K67").Value = 3079.8
$ws.Range("L67").Value = 3080
$ws.Range("M67").Value = -2221.8
$ws.Range("N67").Value = -4796
# hunk @ diff line 5554
$ws.Range("H98").Value = 692.08
$ws.Range("I98").Value = 521.44446
$ws.Range("J98").Value = 1130.8572
$ws.Range("K98").Value = 521.44446
$ws.Range("L98").Value = 1130.8572
$ws.Range("M98").Value = 976.55554
$ws.Range("N98").Value = -4126.8572
# hunk @ diff line 6751
$ws.Range("H122").Value = 692.08
$ws.Range("I122").Value = 521.44446
$ws.Range("J122").Value = 1130.8572
$ws.Range("K122").Value = 1564.33338
$ws.Range("L122").Value = 3392.5716
$ws.Range("M122").Value = 885.66662
$ws.Range("N122").Value = -8292.571599999999
# hunk @ diff line 7357
$ws.Range("H134").Value = 61555.715
$ws.Range("J134").Value = 61555.715
$ws.Range("L134").Value = 61555.715
$ws.Range("N134").Value = -71695.715
# hunk @ diff line 7406
$ws.Range("H135").Value = 808.2
$ws.Range("I135").Value = 457.6
$ws.Range("J135").Value = 2210.6
$ws.Range("K135").Value = 4118.400000000001
$ws.Range("L135").Value = 19895.4
$ws.Range("M135").Value = -1583.400000000001
$ws.Range("N135").Value = -24965.4
# hunk @ diff line 7559
$ws.Range("H138").Value = 989.47
$ws.Range("I138").Value = 528.3889
$ws.Range("J138").Value = 1530.7391
$ws.Range("K138").Value = 1585.1667
$ws.Range("L138").Value = 4592.2173
$ws.Range("M138").Value = 3554.8333
$ws.Range("N138").Value = -14872.2173
# hunk @ diff line 7660
$ws.Range("H140").Value = 50585
$ws.Range("J140").Value = 50585
$ws.Range("L140").Value = 50585
$ws.Range("N140").Value = -60945

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# hunk @ diff line 10859
$ws.Range("H63").Value = 23392.318
$ws.Range("I63").Value = 77831.5
$ws.Range("J63").Value = 4727.457
$ws.Range("K63").Value = 77831.5
$ws.Range("L63").Value = 4727.457
$ws.Range("M63").Value = -77145.5
$ws.Range("N63").Value = -6099.457
# hunk @ diff line 11009
$ws.Range("H66").Value = 23392.318
$ws.Range("I66").Value = 77831.5
$ws.Range("J66").Value = 4727.457
$ws.Range("K66").Value = 389157.5
$ws.Range("L66").Value = 23637.285
$ws.Range("M66").Value = -385725.5
$ws.Range("N66").Value = -30501.285

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# hunk @ diff line 15804
$ws.Range("H22").Value = 255.54286
$ws.Range("I22").Value = 251.375
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 251.375
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -78.375
$ws.Range("N22").Value = -646
# hunk @ diff line 18625
$ws.Range("H80").Value = 1281.4
$ws.Range("I80").Value = 1659.5385
$ws.Range("J80").Value = 579.1429000000001
$ws.Range("K80").Value = 1659.5385
$ws.Range("L80").Value = 579.1429000000001
$ws.Range("M80").Value = -661.5385000000001
$ws.Range("N80").Value = -2575.1429
# hunk @ diff line 18778
$ws.Range("H83").Value = 1281.4
$ws.Range("I83").Value = 1659.5385
$ws.Range("J83").Value = 579.1429000000001
$ws.Range("K83").Value = 8297.692500000001
$ws.Range("L83").Value = 2895.7145
$ws.Range("M83").Value = -3305.692500000001
$ws.Range("N83").Value = -12879.7145

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# hunk @ diff line 22575
$ws.Range("H19").Value = 1312.5
$ws.Range("I19").Value = 1312.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1312.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1142.5
$ws.Range("N19").ClearContents()
# hunk @ diff line 22826
$ws.Range("H24").Value = 1312.5
$ws.Range("I24").Value = 1312.5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1312.5
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -1142.5
$ws.Range("N24").ClearContents()
# hunk @ diff line 23169
$ws.Range("H31").Value = 4401.825
$ws.Range("I31").Value = 1325.0303
$ws.Range("J31").Value = 7786.3
$ws.Range("K31").Value = 1325.0303
$ws.Range("L31").Value = 7786.3
$ws.Range("M31").Value = -1030.0303
$ws.Range("N31").Value = -8376.299999999999
# hunk @ diff line 23316
$ws.Range("H34").Value = 4401.825
$ws.Range("I34").Value = 1325.0303
$ws.Range("J34").Value = 7786.3
$ws.Range("K34").Value = 1325.0303
$ws.Range("L34").Value = 7786.3
$ws.Range("M34").Value = -1123.0303
$ws.Range("N34").Value = -8190.3
# hunk @ diff line 28103
$ws.Range("H132").Value = 3206581
$ws.Range("I132").Value = 1239.2162
$ws.Range("J132").Value = 11113090
$ws.Range("K132").Value = 3717.6486
$ws.Range("L132").Value = 33339270
$ws.Range("M132").Value = -1187.6486
$ws.Range("N132").Value = -33344330

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# hunk @ diff line 28849
$ws.Range("H5").Value = 657.26666
$ws.Range("I5").Value = 346.58334
$ws.Range("J5").Value = 1900
$ws.Range("K5").Value = 1039.75002
$ws.Range("L5").Value = 5700
$ws.Range("M5").Value = -927.7500199999999
$ws.Range("N5").Value = -5924
# hunk @ diff line 34387
$ws.Range("H113").Value = 481.07272
$ws.Range("I113").Value = 464.89474
$ws.Range("J113").Value = 489.6111
$ws.Range("K113").Value = 1394.68422
$ws.Range("L113").Value = 1468.8333
$ws.Range("M113").Value = 775.3157799999999
$ws.Range("N113").Value = -5808.8333
# hunk @ diff line 35051
$ws.Range("H126").Value = 2872.5
$ws.Range("I126").Value = 1230
$ws.Range("J126").Value = 2922.2727
$ws.Range("K126").Value = 3690
$ws.Range("L126").Value = 8766.8181
$ws.Range("M126").Value = 1250
$ws.Range("N126").Value = -18646.8181
# hunk @ diff line 35513
$ws.Range("H135").Value = 657.26666
$ws.Range("I135").Value = 346.58334
$ws.Range("J135").Value = 1900
$ws.Range("K135").Value = 3119.25006
$ws.Range("L135").Value = 17100
$ws.Range("M135").Value = -584.2500600000003
$ws.Range("N135").Value = -22170

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# hunk @ diff line 36317
$ws.Range("H9").Value = 8320
$ws.Range("I9").Value = 400
$ws.Range("J9").Value = 40000
$ws.Range("K9").Value = 400
$ws.Range("L9").Value = 40000
$ws.Range("M9").Value = -230
$ws.Range("N9").Value = -40340
# hunk @ diff line 42302
$ws.Range("H132").Value = 2602.8125
$ws.Range("I132").Value = 2069.25
$ws.Range("J132").Value = 3669.9375
$ws.Range("K132").Value = 6207.75
$ws.Range("L132").Value = 11009.8125
$ws.Range("M132").Value = -3677.75
$ws.Range("N132").Value = -16069.8125
# hunk @ diff line 42752
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# hunk @ diff line 55877
$ws.Range("H126").Value = 1178.3889
$ws.Range("I126").Value = 1200.0667
$ws.Range("J126").Value = 1070
$ws.Range("K126").Value = 3600.2001
$ws.Range("L126").Value = 3210
$ws.Range("M126").Value = -1130.2001
$ws.Range("N126").Value = -8150
# hunk @ diff line 56321
$ws.Range("H135").Value = 100594.164
$ws.Range("J135").Value = 100594.164
$ws.Range("L135").Value = 100594.164
$ws.Range("N135").Value = -110734.164
# hunk @ diff line 56370
$ws.Range("H136").Value = 1945.8152
$ws.Range("I136").Value = 1911.1666
$ws.Range("J136").Value = 2010.7812
$ws.Range("K136").Value = 6207.75
$ws.Range("L136").Value = 6032.3436
$ws.Range("M136").Value = -3183.4998
$ws.Range("N136").Value = -11132.3436
# hunk @ diff line 56618
$ws.Range("H141").Value = 68071.664
$ws.Range("J141").Value = 68071.664
$ws.Range("L141").Value = 68071.664
$ws.Range("N141").Value = -78431.664
